# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates DAMSLTag (col I) and DialogAct (col J)
# values for the rows affected by the re-annotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5;   Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 12;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 39;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 41;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 56;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 76;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 79;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 80;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 82;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 85;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 87;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 88;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 90;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 103; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 109; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 114; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 117; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 129; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 130; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 142; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 149; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 151; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 158; Tag = "sv"; Act = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
